# Apply changes described by the commit "Code Changes for Write Excel and serverRequestType"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the XPath strings in column B (rows 27-32) ---
# Remove the "/YOUI/YOUI" segment from the Comparer XPath expressions.
$ws.Range("B27").Value = "session/data/policy/line/Comparer/12MonthsPremium/Premium"
$ws.Range("B28").Value = "session/data/policy/line/Comparer/12MonthsPremium/GST"
$ws.Range("B29").Value = "session/data/policy/line/Comparer/12MonthsPremium/TotalPremium"
$ws.Range("B30").Value = "session/data/policy/line/Comparer/6MonthsPremium/Premium"
$ws.Range("B31").Value = "session/data/policy/line/Comparer/6MonthsPremium/GST"
$ws.Range("B32").Value = "session/data/policy/line/Comparer/6MonthsPremium/TotalPremium"

# --- Row 4 (EngineCapacity): mark InputOutputType as "Input" ---
$ws.Range("C4").Value = "Input"

# --- Row 16 (Make) and Row 17 (Model): clear the InputOutputType value ---
$ws.Range("C16").Value = ""
$ws.Range("C17").Value = ""

# --- Update the view: current selection ---
$ws.Range("B27").Select()
